# The deck ships two theme parts:
#   ppt/theme/theme1.xml  ("Office Theme" colours) -> used by the Notes Master
#   ppt/theme/theme2.xml  ("Integral" colours)      -> used by the Slide Master
#     (and referenced from the presentation's own "theme" relationship)
#
# The commit swaps the *content* of the two theme parts: the Slide Master
# (and therefore every slide) switches from the "Integral" colour scheme to
# the default "Office Theme" colour scheme, while the Notes Master swaps the
# other way around.
#
# PowerPoint's object model doesn't expose the raw theme XML, but it does let
# us repaint every slot of the active theme's colour scheme via
# Slide.ThemeColorScheme(index).RGB, which is exactly the mechanism PowerPoint
# itself uses under the hood when a different theme/colour-scheme is applied
# from the Design tab. Re-point every one of the 12 theme colour slots on the
# Slide Master's theme to the stock "Office Theme" values.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> target "Office Theme" colour (hex RRGGBB), in the scheme's native
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
$officeTheme = @(
    @(0x00, 0x00, 0x00), # 1  dk1
    @(0xFF, 0xFF, 0xFF), # 2  lt1
    @(0x44, 0x54, 0x6A), # 3  dk2
    @(0xE7, 0xE6, 0xE6), # 4  lt2
    @(0x5B, 0x9B, 0xD5), # 5  accent1
    @(0xED, 0x7D, 0x31), # 6  accent2
    @(0xA5, 0xA5, 0xA5), # 7  accent3
    @(0xFF, 0xC0, 0x00), # 8  accent4
    @(0x44, 0x72, 0xC4), # 9  accent5
    @(0x70, 0xAD, 0x47), # 10 accent6
    @(0x05, 0x63, 0xC1), # 11 hlink
    @(0x95, 0x4F, 0x72)  # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeTheme[$i - 1]
    $tcs.Item($i).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
